# Applies the weekly Fruta/Hortaliza refresh: each data row (3-22) keeps its
# market/category context but receives a new observation date plus updated
# quality/volume/price/origin figures, matching the published diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    3 = @{ 'D' = 44217; 'I' = 'Extra'; 'J' = 400; 'K' = 2500; 'L' = 2500; 'M' = 2500; 'N' = '$/unidad'; 'O' = 'Región de O''Higgins'; 'P' = 2500 }
    4 = @{ 'D' = 44217; 'J' = 280; 'K' = 2000; 'L' = 2000; 'M' = 2000; 'N' = '$/unidad'; 'O' = 'Región de O''Higgins'; 'P' = 2000 }
    5 = @{ 'D' = 44497; 'J' = 250 }
    6 = @{ 'D' = 44194; 'I' = 'Extra'; 'J' = 120; 'K' = 3500; 'L' = 3500; 'M' = 3500; 'N' = '$/unidad'; 'O' = 'Región de O''Higgins'; 'P' = 3500 }
    7 = @{ 'D' = 44194; 'J' = 200; 'K' = 3000; 'L' = 3000; 'M' = 3000; 'O' = 'Región de O''Higgins'; 'P' = 3000 }
    8 = @{ 'D' = 44495; 'J' = 200 }
    9 = @{ 'D' = 44504; 'I' = 'Primera'; 'J' = 200; 'K' = 800; 'L' = 800; 'M' = 800; 'N' = '$/kilo (volumen en unidades)'; 'O' = 'Perú'; 'P' = 800 }
    10 = @{ 'D' = 44510; 'J' = 250; 'K' = 800; 'L' = 800; 'M' = 800; 'N' = '$/kilo (volumen en unidades)'; 'O' = 'Perú'; 'P' = 800 }
    11 = @{ 'D' = 44483; 'I' = 'Primera'; 'J' = 120; 'K' = 800; 'L' = 800; 'M' = 800; 'N' = '$/kilo (volumen en unidades)'; 'O' = 'Perú'; 'P' = 800 }
    12 = @{ 'D' = 44167; 'J' = 400; 'K' = 5000; 'L' = 5000; 'M' = 5000; 'P' = 5000 }
    13 = @{ 'D' = 44167; 'H' = 'Sin especificar'; 'I' = 'Segunda'; 'J' = 560; 'K' = 3000; 'L' = 3000; 'M' = 3000; 'P' = 3000 }
    14 = @{ 'D' = 44167; 'H' = 'Sin especificar'; 'I' = 'Tercera'; 'J' = 450 }
    15 = @{ 'D' = 44491; 'H' = 'Sin especificar'; 'I' = 'Primera'; 'J' = 150; 'K' = 800; 'L' = 800; 'M' = 800; 'N' = '$/kilo (volumen en unidades)'; 'O' = 'Perú'; 'P' = 800 }
    16 = @{ 'D' = 44477; 'H' = 'Sin especificar'; 'I' = 'Primera'; 'J' = 80; 'K' = 800; 'L' = 800; 'M' = 800; 'N' = '$/kilo (volumen en unidades)'; 'O' = 'Perú'; 'P' = 800 }
    17 = @{ 'D' = 44305; 'J' = 100; 'K' = 2500; 'L' = 2500; 'M' = 2500; 'N' = '$/unidad'; 'P' = 2500 }
    18 = @{ 'D' = 44223; 'H' = 'Americana O Klondike'; 'I' = 'Extra'; 'J' = 340; 'K' = 2500; 'L' = 2500; 'M' = 2500; 'N' = '$/unidad'; 'O' = 'Región de O''Higgins'; 'P' = 2500 }
    19 = @{ 'D' = 44223; 'H' = 'Americana O Klondike'; 'J' = 400; 'K' = 2000; 'L' = 2000; 'M' = 2000; 'O' = 'Región de O''Higgins'; 'P' = 2000 }
    20 = @{ 'D' = 44223; 'H' = 'Americana O Klondike'; 'I' = 'Segunda'; 'J' = 300; 'K' = 1500; 'L' = 1500; 'M' = 1500; 'P' = 1500 }
    21 = @{ 'D' = 44223; 'H' = 'Americana O Klondike'; 'I' = 'Tercera'; 'J' = 160; 'K' = 1000; 'L' = 1000; 'M' = 1000; 'P' = 1000 }
    22 = @{ 'D' = 44312; 'I' = 'Primera'; 'J' = 180; 'K' = 2500; 'L' = 2500; 'M' = 2500; 'O' = 'Perú'; 'P' = 2500 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
